# ---------------------------------------------------------------------------
# This script reproduces (against the Word object model) the changes described
# by the commit: two paragraphs of exercise text are tidied up / extended in
# "CiclosAbstraccion1Basicos.docx".
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: the paragraph "Utilizando ciclos, ... Por ejemplo para A = 0 ..."
# had the word "ejemplo" wrapped in its own run flanked by proofErr
# (grammar-check) markers. Re-typing the sentence as a plain Find & Replace
# collapses it back down to a single run and drops the now-stale proofErr
# markers, while leaving the visible text untouched.
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute(
    "Por ejemplo para A = 0 y   B = 5, imprimir" + [char]0x00ED + "a:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Por ejemplo para A = 0 y   B = 5, imprimir" + [char]0x00ED + "a:",
    2
) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: locate the MCD exercise paragraph ("Realice una funcion que
# calcule el maximo comun divisor (MCD) ..."). It gets reflowed (the
# misspelled "aquel" no longer needs the gramStart/gramEnd proofErr wrapper)
# and a brand new paragraph describing a "is this number prime?" exercise is
# appended right after it, in the same numbered-list style. The hidden
# "_GoBack" bookmark (Word's "last edit" marker), which used to sit right
# before the MCD text, now sits at the very end of the new paragraph -- i.e.
# at the point of the most recent edit.
# ---------------------------------------------------------------------------
$find2 = $d.Content.Find
$find2.ClearFormatting()
$found2 = $find2.Execute("MCD", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not locate the MCD exercise paragraph"
}
$mcdPara = $find2.Parent.Paragraphs(1)

$xmlFrag2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Realice una función que calcule el máximo común divisor (MCD) entre 2 números enteros positivos. Recordemos que el MCD de 2 números, es aqu</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>e</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>l número más grande, que logra dividir a ambos números a la vez, por ejemplo el MCD entre 12 y 8, sería el número 4, pues es el número más grande que logra dividir de manera exacta al 8 y al 12 a la vez.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Realice un</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">a función que indique si un numero es primo o no. Luego use esa función para </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>m</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>o</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>str</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>ar</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> los números primos desde 1 hasta N, dónde N será un valor digitado por el usuario</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$mcdPara.Range.InsertXML($xmlFrag2)
